$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17 (pushing the old blank spacer row 17,
# and everything below it, down by one). The old row 17 becomes row 18,
# which is where the new "DoSleep" wait step goes - this mirrors the
# formatting seen in the target file, where the new spacer row keeps its
# own distinct style and the step below reuses the old spacer's style.
$ws.Rows.Item(17).Insert()
$ws.Range("A17").NumberFormat = "General"

# Populate the newly freed-up row 18 with the new Global DoSleep action -
# an extra wait needed to support the latest D365 version.
$ws.Range("B18").Value = "Action"
$ws.Range("C18").Value = "Global"
$ws.Range("D18").Value = "DoSleep"
$ws.Range("E18").Value = "millis"
$ws.Range("F18").Value = "number"

# Param Value is stored as text throughout this sheet (even numeric-looking
# values like "5000"), so force a text number format before assigning it.
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "5000"
